# Apply the "Updated cryptos list" data refresh (Fri Jun 23 08:28:34 UTC 2023).
#
# For every changed cell we just set .Value to the new text coming from the
# source feed. For the Price column (D) the new text is often something that
# *looks* like a plain number ("9.100", "103.10", "0.000007841", ...). A bare
# `.Value = "..."` assignment there would let Excel's type-inference turn it
# into a real number and silently drop the significant trailing/leading zeros
# (e.g. "9.100" -> 9.1), which would not match the source data (plain text).
# To avoid that, those cells are first switched to an explicit text format,
# written, and then reset back to the workbook's normal (unstyled) cell style
# so no stray number formatting is left behind on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.963.77"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").Value = "1.881.45"
$ws.Range("E3").Value = "  -1.42%  "
$cell = $ws.Range("D4")
$cell.NumberFormat = "@"
$cell.Value = "0.9991"
$cell.Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = "243.34"
$cell.Style = "Normal"
$ws.Range("E5").Value = "  -3.40%  "
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = "0.9985"
$cell.Style = "Normal"
$ws.Range("E6").Value = "  -0.14%  "
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = "0.4932"
$cell.Style = "Normal"
$ws.Range("E7").Value = "  -3.65%  "
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = "0.2946"
$cell.Style = "Normal"
$ws.Range("E8").Value = "  -2.02%  "
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = "0.06645"
$cell.Style = "Normal"
$ws.Range("E9").Value = "  -2.44%  "
$ws.Range("D10").Value = "1.875.16"
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = "16.75"
$cell.Style = "Normal"
$ws.Range("E11").Value = "  -3.25%  "
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = "0.07209"
$cell.Style = "Normal"
$ws.Range("E12").Value = "  -1.67%  "
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = "0.6687"
$cell.Style = "Normal"
$ws.Range("E13").Value = "  -4.55%  "
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = "86.42"
$cell.Style = "Normal"
$ws.Range("E14").Value = "  -0.44%  "
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = "4.886"
$cell.Style = "Normal"
$ws.Range("E15").Value = "  -0.62%  "
$ws.Range("D16").Value = "29.931.10"
$ws.Range("E16").Value = "  -0.51%  "
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = "0.000007841"
$cell.Style = "Normal"
$ws.Range("E17").Value = "  -4.71%  "
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = "12.82"
$cell.Style = "Normal"
$ws.Range("E19").Value = "  -2.03%  "
$ws.Range("D20").Value = "2.119.00"
$ws.Range("E20").Value = "  -1.68%  "
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = "0.9989"
$cell.Style = "Normal"
$ws.Range("E21").Value = "  -0.03%  "
$cell = $ws.Range("D22")
$cell.NumberFormat = "@"
$cell.Value = "4.791"
$cell.Style = "Normal"
$ws.Range("E22").Value = "  -0.97%  "
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = "5.881"
$cell.Style = "Normal"
$ws.Range("E23").Value = "  +2.03%  "
$cell = $ws.Range("D24")
$cell.NumberFormat = "@"
$cell.Value = "9.100"
$cell.Style = "Normal"
$ws.Range("E24").Value = "  -1.78%  "
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = "150.29"
$cell.Style = "Normal"
$ws.Range("E25").Value = "  +1.48%  "
$cell = $ws.Range("D26")
$cell.NumberFormat = "@"
$cell.Value = "142.86"
$cell.Style = "Normal"
$ws.Range("E26").Value = "  +5.71%  "
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = "17.07"
$cell.Style = "Normal"
$ws.Range("E27").Value = "  -0.04%  "
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = "1.924"
$cell.Style = "Normal"
$ws.Range("E28").Value = "  -3.94%  "
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = "1.387"
$cell.Style = "Normal"
$ws.Range("E29").Value = "  -1.03%  "
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = "4.219"
$cell.Style = "Normal"
$ws.Range("E30").Value = "  -1.34%  "
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = "0.08786"
$cell.Style = "Normal"
$ws.Range("E31").Value = "  -0.47%  "
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = "3.995"
$cell.Style = "Normal"
$ws.Range("E32").Value = "  -0.31%  "
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = "0.05055"
$cell.Style = "Normal"
$ws.Range("E33").Value = "  -0.15%  "
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = "0.7145"
$cell.Style = "Normal"
$ws.Range("E34").Value = "  -0.65%  "
$cell = $ws.Range("D35")
$cell.NumberFormat = "@"
$cell.Value = "1.116"
$cell.Style = "Normal"
$ws.Range("E35").Value = "  -2.36%  "
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = "2.667"
$cell.Style = "Normal"
$ws.Range("E36").Value = "  -0.85%  "
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = "0.01796"
$cell.Style = "Normal"
$ws.Range("E37").Value = "  +5.66%  "
$cell = $ws.Range("D38")
$cell.NumberFormat = "@"
$cell.Value = "2.701"
$cell.Style = "Normal"
$ws.Range("E38").Value = "  -3.85%  "
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = "2.175"
$cell.Style = "Normal"
$ws.Range("E39").Value = "  -3.89%  "
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = "0.9317"
$cell.Style = "Normal"
$ws.Range("E40").Value = "  -3.24%  "
$ws.Range("B41").Value = "FraxShare"
$ws.Range("C41").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = "5.772"
$cell.Style = "Normal"
$ws.Range("E41").Value = "  -6.43%  "
$ws.Range("B42").Value = "TheSandbox"
$ws.Range("C42").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = "0.4236"
$cell.Style = "Normal"
$ws.Range("E42").Value = "  -1.86%  "
$ws.Range("B43").Value = "PaxDollar"
$ws.Range("C43").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = "0.9985"
$cell.Style = "Normal"
$ws.Range("E43").Value = "  -0.04%  "
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = "103.10"
$cell.Style = "Normal"
$ws.Range("E44").Value = "  -1.75%  "
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = "7.424"
$cell.Style = "Normal"
$ws.Range("E45").Value = "  -2.59%  "
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = "0.1272"
$cell.Style = "Normal"
$ws.Range("E46").Value = "  -0.74%  "
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = "0.05665"
$cell.Style = "Normal"
$ws.Range("E47").Value = "  -1.31%  "
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = "32.59"
$cell.Style = "Normal"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = "8.296"
$cell.Style = "Normal"
$ws.Range("E49").Value = "  -1.92%  "
$ws.Range("B50").Value = "Decentraland"
$ws.Range("C50").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = "0.3772"
$cell.Style = "Normal"
$ws.Range("E50").Value = "  -1.27%  "
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = "56.12"
$cell.Style = "Normal"
$ws.Range("E51").Value = "  -1.31%  "
